$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values: B2 226 -> 243, B3 171 -> 226
$ws.Range("B2").Value = 243
$ws.Range("B3").Value = 226

# Remove row 4 entirely (A4=2, B4=72) since the data now only spans A1:B3
$ws.Range("A4:B4").Delete()
